$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "SACRAMENTO AREA TOTALS" text from B2 to A2, and set B2 to "Totals"
$ws.Range("A2").Value2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = "Totals"

# Match column A's width to column B's width, and select B4
$ws.Columns("A").ColumnWidth = $ws.Columns("B").ColumnWidth
$ws.Range("B4").Select()
